$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the old "Delete Account" step rows (19-20) that are being replaced
# by the new test-scenario rows below.
$ws.Range("F19:H20").ClearContents()

# TestScenario_1 / TestCase_1 - New Account
$ws.Range("A2").Value = "TestScenario_1"
$ws.Range("B2").Value = "TestScenario_1.TestCase_1"
$ws.Range("C2").Value = "New Account"
$ws.Range("D2").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F2").Value = "Step 1"
$ws.Range("G2").Value = "Click on the Account tab, and click on New button"
$ws.Range("H2").Value = "User should be navigated to the New  Account Page"
$ws.Range("I2").Value = "Approved"

# TestScenario_2 / TestCase_1 - View Account
$ws.Range("A9").Value = "TestScenario_2"
$ws.Range("B9").Value = "TestScenario_2.TestCase_1"
$ws.Range("C9").Value = "View Account"
$ws.Range("D9").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F9").Value = "Step 1"
$ws.Range("G9").Value = "Click on the Account tab,  and select a Account "
$ws.Range("H9").Value = "User should be navigated to the Account Page"
$ws.Range("I9").Value = "Approved"

# TestScenario_3 / TestCase_1 - Edit Account
$ws.Range("A11").Value = "TestScenario_3"
$ws.Range("B11").Value = "TestScenario_3.TestCase_1"
$ws.Range("C11").Value = "Edit Account"
$ws.Range("D11").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F11").Value = "Step 1"
$ws.Range("G11").Value = "Click on the Account tab,  and click on existing Account to modify"
$ws.Range("H11").Value = "User is navigated to the Account Details page"
$ws.Range("I11").Value = "Approved"

# TestScenario_4 / TestCase_1 - Delete Account
$ws.Range("A18").Value = "TestScenario_4"
$ws.Range("B18").Value = "TestScenario_4.TestCase_1"
$ws.Range("C18").Value = "Delete Account"
$ws.Range("D18").Value = "User Needs to Login to Salesforce, from the browser with correct credentials"
$ws.Range("F18").Value = "Step 1"
$ws.Range("G18").Value = "Click on the Account tab,  and select the existing  Account to delete"
$ws.Range("H18").Value = "User is navigated to the Account Details page"
$ws.Range("I18").Value = "Approved"

# Match the author's final selection/scroll state
$ws.Range("I18").Select()
